# resdoc-content.docx "minor tweaks" edit
# 1. Heading3 paragraph style: language fr-CA -> en-CA (also cascades to the
#    Heading3Char linked character style).
# 2. Compact paragraph style: drop the explicit 10pt override so the font
#    size resolves to the inherited 11pt value from BodyText/Normal.
# 3. References table: widen the two gridCol columns.

$d = $word.ActiveDocument

# --- 1. Heading 3 / Heading 3 Char styles: fr-CA -> en-CA -------------------
$heading3 = $d.Styles.Item("Heading3")
$heading3.LanguageID = "en-CA"

$heading3Char = $d.Styles.Item("Heading3Char")
$heading3Char.LanguageID = "en-CA"

# --- 2. Compact style: remove the explicit 10pt size override --------------
$compact = $d.Styles.Item("Compact")
$compact.Font.Size = 11

# --- 3. References table column widths --------------------------------
$tbl = $d.Tables(1)
$tbl.Columns(1).Width = 38.35
$tbl.Columns(2).Width = 41.4
